$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("РБ №1")

$xlPasteFormats = -4122

# --- Sheet1: add "Количество" (D) and "Итог:" (E) header columns ---
$ws.Range("D1").Value = "Количество"
$ws.Range("E1").Value = "Итог:"
# E1 is a brand-new cell -- give it the same boxed header style as D1/C1/etc.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial($xlPasteFormats)

# Quantities for each product row (2-24), aligned to existing rows
$qty = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 1
    6 = 1
    7 = 30
    8 = 1
    9 = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 5
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 5
    24 = 20
}

$grandTotal = 0
for ($row = 2; $row -le 24; $row++) {
    $q = $qty[$row]
    $price = $ws.Cells.Item($row, 3).Value2
    $lineTotal = $price * $q
    $ws.Cells.Item($row, 4).Value = $q
    $ws.Cells.Item($row, 5).Value = $lineTotal
    $grandTotal = $grandTotal + $lineTotal
}

# Row 24's quantity cell (D24) used to carry the boxed note style + text;
# that note moves to a brand-new F24 cell, so D24 reverts to a plain cell.
$ws.Range("C24").Copy()
$ws.Range("D24").PasteSpecial($xlPasteFormats)

$ws.Range("F24").Value = "Сумма / Процент РБ"
$ws.Range("A25").Copy()
$ws.Range("F24").PasteSpecial($xlPasteFormats)

# Unmerge A25:B25, and clear A25 / C25 (keep their box style)
$ws.Range("A25:B25").UnMerge()
$ws.Range("A25").Value = $null
$ws.Range("C25").Value = $null

# New summary row 25: D25 = "Итог:" label, E25 = grand total, F25 = 5000
$ws.Range("D25").Value = "Итог:"
$ws.Range("B25").Copy()
$ws.Range("D25").PasteSpecial($xlPasteFormats)

$ws.Range("E25").Value = $grandTotal
$ws.Range("B25").Copy()
$ws.Range("E25").PasteSpecial($xlPasteFormats)

$ws.Range("F25").Value = 5000
$ws.Range("B25").Copy()
$ws.Range("F25").PasteSpecial($xlPasteFormats)

# --- Sheet2 ("РБ №1"): update amounts ---
$ws2.Range("D2").Value = 5000
$ws2.Range("E2").Value = 5000
$ws2.Range("E3").Value = 5000
